# buglist&newfeature.xlsx update
# 1 meetplayer implement dlna push_to_dmr
# 2 win32 testDlg support vlc play

$wb = $excel.ActiveWorkbook
$wsBug = $wb.Worksheets.Item("bug")
$wsNew = $wb.Worksheets.Item("newfeature")

# ---------------------------------------------------------------------------
# "bug" sheet - new rows 29 (#28), 30 (#29), 31 (#30)
# ---------------------------------------------------------------------------

# Row 29 (bug #28) - keeps the plain (non-highlighted) style already on the row
$wsBug.Range("C29").Value = "N/A"
$wsBug.Range("D29").Value = "IOS"
$wsBug.Range("E29").Value = 20150129
$wsBug.Range("F29").Value = "由于编译问题vc1解码模块disabled，导致vc1影片无法播放"
$wsBug.Range("G29").Value = "TBD"

# Row 30 (bug #29) - becomes highlighted like row 28, so copy that row's format first
$wsBug.Range("B28:H28").Copy()
$wsBug.Range("B30:H30").PasteSpecial(-4122)
$wsBug.Range("B30").Value = 29
$wsBug.Range("C30").Value = "N/A"
$wsBug.Range("D30").Value = "android"
$wsBug.Range("E30").Value = 20150129
$wsBug.Range("F30").Value = "自有播放器播放 茜拉-想你的夜.WAV pos显示不对"
$wsBug.Range("G30").Value = "tracking"
$wsBug.Range("H30").Value = ""

# Row 31 (bug #30) - also highlighted
$wsBug.Range("B28:H28").Copy()
$wsBug.Range("B31:H31").PasteSpecial(-4122)
$wsBug.Range("B31").Value = 30
$wsBug.Range("C31").Value = "N/A"
$wsBug.Range("D31").Value = "android"
$wsBug.Range("E31").Value = 20150129
$wsBug.Range("F31").Value = "自有播放器播放 陈慧娴-飘雪.ape文件进度条走的很快，声音不对"
$wsBug.Range("G31").Value = ""
$wsBug.Range("H31").Value = ""

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# "newfeature" sheet - new row 23 (#22)
# ---------------------------------------------------------------------------

$wsNew.Range("C23").Value = "meetsdk"
$wsNew.Range("D23").Value = "android"
$wsNew.Range("E23").Value = "扫描本地文件写入db"
$wsNew.Range("F23").Value = "TBD"

# ---------------------------------------------------------------------------
# Selection / active-tab bookkeeping: the new entries were made on "bug", but
# the workbook is left with "newfeature" as the active/selected tab.
# ---------------------------------------------------------------------------

$wsBug.Range("E23").Select()
$wsNew.Activate()
$wsNew.Range("F23").Select()
